$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..98 down to 8..99
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new data point:
# Dia=6, total_venda=16533.99, Mes=8, Ano=2025, Periodo="08/2025"
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 16533.99
$ws.Cells.Item(7, 3).Value = 8
$ws.Cells.Item(7, 4).Value = 2025
$ws.Cells.Item(7, 5).Value = "08/2025"
